# Apply "Append: 2025-09-05 06:26 JST" update to the Lancers sheet.
#
# The scraped-data sheet ("ランサーズ") is refreshed: every surviving row gets a
# new scrape timestamp, new scraped field values, and the six oldest rows
# (14-19) are dropped entirely (the used range shrinks from A1:H19 to A1:H13).
# A handful of the remaining rows also lose their "skill summary" (H) value.
# Column widths for B, D and H are also adjusted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the six rows that fall out of the refreshed scrape (old rows 14-19)
# ---------------------------------------------------------------------------
$ws.Rows("14:19").Delete()

# ---------------------------------------------------------------------------
# 2) Column width changes (values are COM "character" widths; the engine adds
#    a constant 5/6 px-padding when it serialises back to the <col width=.../>
#    attribute, so we compensate by subtracting 5/6 before assigning).
# ---------------------------------------------------------------------------
$pad = 5.0 / 6.0
$ws.Columns("B").ColumnWidth = 48 - $pad
$ws.Columns("D").ColumnWidth = 30 - $pad
$ws.Columns("H").ColumnWidth = 27 - $pad

# ---------------------------------------------------------------------------
# 3) New row data (rows 2-13)
# ---------------------------------------------------------------------------
$timestamp = "2025-09-05 06:26:09"

$rows = @(
    @{ Row=2;  B="詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発"; D="300,000 円 ~ 500,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5377709"; G=245; H="🔥Next.js ◆開発,Node.js ◇アプリ" },
    @{ Row=3;  B="<Next.js、バックエンド開発> ガントチャートアプリの改修製造"; D="300,000 円 ~ 500,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5379158"; G=225; H="🔥Next.js ◆開発 ◇アプリ" },
    @{ Row=4;  B=$null; D=$null; F=$null; G=$null; H=$null },
    @{ Row=5;  B="【React/Vue】新規サービス開発に携わるフロントエンドエンジニア募集(フルリモート可)"; D="500,000 円 ~ 1,000,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5387591"; G=190; H="🔥React ◆開発" },
    @{ Row=6;  B="【注目】公式LINEで診断機能を実現するGPT連動開発依頼"; D=$null; F="https://www.lancers.jp/work/detail/5387629"; G=183; H="🔥GPT ◆開発" },
    @{ Row=7;  B="Flutterなどハイブリッドアプリによる業務アプリの開発(スマートウォッチ)"; D="300,000 円 ~ 500,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5379176"; G=100; H=$null },
    @{ Row=8;  B="【急募】スキースノーボードスクール予約サイトの料金修正依頼"; D="50,000 円 ~ 100,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5374405"; G=38; H="◇サイト" },
    @{ Row=9;  B="【SRE / インフラエンジニア募集(基本リモート・金融系案件)】"; D="500,000 円 ~ 1,000,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5371075"; G=25; H="" },
    @{ Row=10; B="限定公開 PR 限定公開の仕事"; D="500,000 円 ~ 1,000,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5385681"; G=25; H="" },
    @{ Row=11; B="注目 PR 超初級・SE育成の技術研修 サブ講師"; D="500,000 円 ~ 1,000,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5385021"; G=25; H="" },
    @{ Row=12; B="【NAS導入】VPN設定とネットワークドライブの構築支援"; D="20,000 円 ~ 50,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5387417"; G=13; H="" },
    @{ Row=13; B="【急募】Excelで自動シート増加と画像トリミングを実現!"; D="10,000 円 ~ 20,000 円 / 固定"; F="https://www.lancers.jp/work/detail/5387258"; G=10; H="" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $timestamp

    if ($item.B -ne $null) { $ws.Cells.Item($r, 2).Value = $item.B }
    if ($item.D -ne $null) { $ws.Cells.Item($r, 4).Value = $item.D }
    if ($item.F -ne $null) { $ws.Cells.Item($r, 6).Value = $item.F }
    if ($item.G -ne $null) { $ws.Cells.Item($r, 7).Value = $item.G }

    if ($item.H -eq $null) {
        # H unchanged for this row - leave as-is
    } elseif ($item.H -eq "") {
        # H cell removed entirely in the refreshed data
        $ws.Cells.Item($r, 8).ClearContents()
    } else {
        $ws.Cells.Item($r, 8).Value = $item.H
    }
}

# ---------------------------------------------------------------------------
# 4) Hyperlinks: targets for F2:F13 changed (and F14:F19 hyperlinks must go
#    away with the deleted rows). Rebuild the whole collection from the new
#    URLs so the relationship targets match the refreshed cell values.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value2)
}
$ws.Range("F2:F13").Style = "Hyperlink"
